# fix: update the correct static files
#
# Replaces the placeholder "dummy test data" rows/values that had been
# left in the shipped ExcelDatabase.xlsx with the intended sample data,
# and removes the stray test rows that were never meant to ship.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# customer sheet: rename placeholder company names to the "Sample
# Company N" convention, and drop the bogus "Dummy" row (row 8).
# ---------------------------------------------------------------------
$customer = $wb.Worksheets.Item("customer")
$customer.Range("B2").Value = "Sample Company 1"
$customer.Range("B3").Value = "Sample Company 2"
$customer.Range("B4").Value = "Sample Company 3"
$customer.Range("B5").Value = "Sample Company 4"
$customer.Range("B6").Value = "Sample Company 5"
$customer.Range("B7").Value = "Sample Company 6"
$customer.Range("A8:F8").EntireRow.Delete()

# ---------------------------------------------------------------------
# project sheet: fix up the num_samples / total for project 1, and
# drop the "Test Data...." row (row 8).
# ---------------------------------------------------------------------
$project = $wb.Worksheets.Item("project")
$project.Range("C2").Value = 3
$project.Range("W2").Value = 260
$project.Range("A8:X8").EntireRow.Delete()

# ---------------------------------------------------------------------
# services sheet: drop the "Test Service" / "Boris" rows (rows 15-16).
# ---------------------------------------------------------------------
$services = $wb.Worksheets.Item("services")
$services.Range("A15:F16").EntireRow.Delete()

# ---------------------------------------------------------------------
# worker sheet: drop the "John Doe" row (row 8).
# ---------------------------------------------------------------------
$worker = $wb.Worksheets.Item("worker")
$worker.Range("A8:E8").EntireRow.Delete()

# ---------------------------------------------------------------------
# workerprojectbridge sheet: correct the last two bridge rows.
# ---------------------------------------------------------------------
$wpb = $wb.Worksheets.Item("workerprojectbridge")
$wpb.Range("A7").Value = 7
$wpb.Range("B7").Value = 1
$wpb.Range("C7").Value = 2
$wpb.Range("A8").Value = 8
$wpb.Range("B8").Value = 6
$wpb.Range("C8").Value = 6

# ---------------------------------------------------------------------
# projectservicesbridge sheet: correct rows 8-9, drop rows 10-13.
# ---------------------------------------------------------------------
$psb = $wb.Worksheets.Item("projectservicesbridge")
$psb.Range("B8").Value = 2
$psb.Range("C8").Value = 200
$psb.Range("E8").Value = 3
$psb.Range("B9").Value = 1
$psb.Range("C9").Value = 10
$psb.Range("E9").Value = 6
$psb.Range("A10:E13").EntireRow.Delete()

# ---------------------------------------------------------------------
# projectbillbridge sheet: drop rows 8-9.
# ---------------------------------------------------------------------
$pbb = $wb.Worksheets.Item("projectbillbridge")
$pbb.Range("A8:C9").EntireRow.Delete()

# ---------------------------------------------------------------------
# bill sheet: fill in the real billing address / extra charges for
# bills 1 and 2, and drop the leftover test rows (4-5).
# ---------------------------------------------------------------------
$bill = $wb.Worksheets.Item("bill")
$bill.Range("C2").Value = "120 University Place, Glasgow G12 8TA"
$bill.Range("D2").Value = "Extra - Annual charge"
$bill.Range("E2").Value = 100
$bill.Range("H2").Value = 490
$bill.Range("D3").Value = "Sample extra 1"
$bill.Range("E3").Value = 10
$bill.Range("F3").Value = "Sample extra 2"
$bill.Range("G3").Value = 20
$bill.Range("H3").Value = 322.5
$bill.Range("A4:I5").EntireRow.Delete()
